$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Turns the gloss text
#     ". anno ngeiss cách imm a chomalnad"
# into
#     ". an- nongeiss cách imm a chomalnad"
# (splitting "anno" into "an" + "-", and moving the "no" onto the front of
# "ngeiss" to make "nongeiss"), while relocating the document's "_GoBack"
# bookmark from the very end of the document to the new split point
# between "-" and the following space.
#
# The target run layout is:
#   <w:r>. an</w:r>
#   <w:r>-</w:r>
#   <w:bookmarkStart w:name="_GoBack"/><w:bookmarkEnd/>
#   <w:r xml:space="preserve"> </w:r>
#   <w:r>nongeiss cách imm a chomalnad</w:r>
#
# The Word engine coalesces adjacent same-formatted runs whenever an edit
# touches a paragraph, UNLESS a bookmark already sits exactly on the
# boundary at the moment of the edit. So: first drop (zero-width, edit
# free) bookmarks at every boundary we want to keep, THEN replace the text
# of each now-isolated segment in a single whole-range assignment (never
# spanning a bookmark/run boundary, and never using InsertBefore/After,
# which would otherwise weld the new text onto a neighbouring run and tag
# it with a spurious xml:space="preserve"). Finally drop the temporary
# bookmarks, leaving only the real "_GoBack" in its new spot.
# ---------------------------------------------------------------------------

# Locate ". anno" robustly (rather than hard-coding character offsets).
$anchor = $d.Content
$anchor.Find.Execute(". anno", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$startAnno = $anchor.Start      # start of ". anno"
$endAnno = $anchor.End          # end of ". anno" == start of the following space run
$splitAn = $startAnno + 4        # boundary between ". an" and "no"

# 1) Drop the permanent bookmark onto the new split point (between the
#    future "-" and the space run). Adding a bookmark whose name already
#    exists elsewhere in the document relocates it, so this simultaneously
#    removes the stale "_GoBack" that currently sits at the very end of
#    the document (2nd hunk of the change).
$d.Bookmarks.Add("_GoBack", $d.Range($endAnno, $endAnno)) | Out-Null

# 2) Drop a temporary bookmark between "an" and "no".
$d.Bookmarks.Add("TMP_SPLIT_A", $d.Range($splitAn, $splitAn)) | Out-Null

# 3) Drop a temporary bookmark between the space and "ngeiss...".
$spaceStart = $endAnno
$ngeissStart = $spaceStart + 1
$d.Bookmarks.Add("TMP_SPLIT_C", $d.Range($ngeissStart, $ngeissStart)) | Out-Null

# 4) Replace "no" (now isolated between TMP_SPLIT_A and _GoBack) with "-".
$d.Range($splitAn, $endAnno).Text = "-"

# Replacing the 2-character "no" with the 1-character "-" shifts every
# later offset back by one.
$shift = -1
$ngeissStart = $ngeissStart + $shift

# 5) Replace "ngeiss cách imm a chomalnad" (now isolated after
#    TMP_SPLIT_C) with "nongeiss cách imm a chomalnad" in one whole-range
#    assignment so the run doesn't pick up a spurious xml:space.
$oldNgeissLen = "ngeiss cách imm a chomalnad".Length
$d.Range($ngeissStart, $ngeissStart + $oldNgeissLen).Text = "nongeiss cách imm a chomalnad"

# 6) Remove the temporary bookmarks; "_GoBack" (added in step 1) is left
#    in place between "-" and the space.
$d.Bookmarks("TMP_SPLIT_A").Delete()
$d.Bookmarks("TMP_SPLIT_C").Delete()

Write-Output ("Final text: " + $d.Range($startAnno, $startAnno + 36).Text)
